$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 11 new rows after the header, shifting existing data down by 11
$ws.Range("2:12").Insert()

# Populate the 11 newly inserted sensor-data rows (rows 2-12)
$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = "falling"
$ws.Cells.Item(2, 3).Value = -1.931140422821045
$ws.Cells.Item(2, 4).Value = 9.274446487426758
$ws.Cells.Item(2, 5).Value = 0.347445011138916
$ws.Cells.Item(2, 6).Value = 0.308792382478714
$ws.Cells.Item(2, 7).Value = -0.2591595947742462
$ws.Cells.Item(2, 8).Value = 0.3060434758663177

$ws.Cells.Item(3, 1).Value = 100
$ws.Cells.Item(3, 2).Value = "falling"
$ws.Cells.Item(3, 3).Value = -2.475548833608627
$ws.Cells.Item(3, 4).Value = 9.384642362594603
$ws.Cells.Item(3, 5).Value = 0.7847917079925537
$ws.Cells.Item(3, 6).Value = 0.197004035115242
$ws.Cells.Item(3, 7).Value = 0.0836885422468185
$ws.Cells.Item(3, 8).Value = 0.0430659987032413

$ws.Cells.Item(4, 1).Value = 200
$ws.Cells.Item(4, 2).Value = "falling"
$ws.Cells.Item(4, 3).Value = -2.621871948242188
$ws.Cells.Item(4, 4).Value = 9.314098954200743
$ws.Cells.Item(4, 5).Value = 1.41282993555069
$ws.Cells.Item(4, 6).Value = 0.0916297882795333
$ws.Cells.Item(4, 7).Value = 0.034972034394741
$ws.Cells.Item(4, 8).Value = 0.0708603709936142

$ws.Cells.Item(5, 1).Value = 300
$ws.Cells.Item(5, 2).Value = "falling"
$ws.Cells.Item(5, 3).Value = -2.868601083755493
$ws.Cells.Item(5, 4).Value = 9.443870902061462
$ws.Cells.Item(5, 5).Value = 0.9555243626236917
$ws.Cells.Item(5, 6).Value = 0.0360410511493682
$ws.Cells.Item(5, 7).Value = 0.06902777403593061
$ws.Cells.Item(5, 8).Value = -0.0134390350431203

$ws.Cells.Item(6, 1).Value = 400
$ws.Cells.Item(6, 2).Value = "falling"
$ws.Cells.Item(6, 3).Value = -2.620113015174865
$ws.Cells.Item(6, 4).Value = 9.546792268753052
$ws.Cells.Item(6, 5).Value = 0.6990440487861633
$ws.Cells.Item(6, 6).Value = -0.0181732401251792
$ws.Cells.Item(6, 7).Value = 0.0655152946710586
$ws.Cells.Item(6, 8).Value = -0.0574213340878486

$ws.Cells.Item(7, 1).Value = 500
$ws.Cells.Item(7, 2).Value = "falling"
$ws.Cells.Item(7, 3).Value = -2.44504114985466
$ws.Cells.Item(7, 4).Value = 9.541788041591644
$ws.Cells.Item(7, 5).Value = 0.372002582065761
$ws.Cells.Item(7, 6).Value = -0.1020144969224929
$ws.Cells.Item(7, 7).Value = -0.0171042270958423
$ws.Cells.Item(7, 8).Value = -0.0577267669141292

$ws.Cells.Item(8, 1).Value = 600
$ws.Cells.Item(8, 2).Value = "falling"
$ws.Cells.Item(8, 3).Value = -2.28140389919281
$ws.Cells.Item(8, 4).Value = 9.524857640266418
$ws.Cells.Item(8, 5).Value = -0.01871592737734312
$ws.Cells.Item(8, 6).Value = -0.0662788823246955
$ws.Cells.Item(8, 7).Value = -0.0591012127697467
$ws.Cells.Item(8, 8).Value = -0.0519235469400882

$ws.Cells.Item(9, 1).Value = 700
$ws.Cells.Item(9, 2).Value = "falling"
$ws.Cells.Item(9, 3).Value = -2.462455779314041
$ws.Cells.Item(9, 4).Value = 9.538427114486694
$ws.Cells.Item(9, 5).Value = -0.03734804317355161
$ws.Cells.Item(9, 6).Value = -0.0091629782691597
$ws.Cells.Item(9, 7).Value = 0.016951510682702
$ws.Cells.Item(9, 8).Value = 0.0032070425804704

$ws.Cells.Item(10, 1).Value = 800
$ws.Cells.Item(10, 2).Value = "falling"
$ws.Cells.Item(10, 3).Value = -2.66546654701233
$ws.Cells.Item(10, 4).Value = 9.521270275115967
$ws.Cells.Item(10, 5).Value = 0.1833332777023315
$ws.Cells.Item(10, 6).Value = -0.0189368221908807
$ws.Cells.Item(10, 7).Value = 0.0145080499351024
$ws.Cells.Item(10, 8).Value = 0.0222965814173221

$ws.Cells.Item(11, 1).Value = 900
$ws.Cells.Item(11, 2).Value = "falling"
$ws.Cells.Item(11, 3).Value = -2.735388696193695
$ws.Cells.Item(11, 4).Value = 9.501047194004059
$ws.Cells.Item(11, 5).Value = 0.2642159881070256
$ws.Cells.Item(11, 6).Value = -0.0164933614432811
$ws.Cells.Item(11, 7).Value = -0.0006108652451075
$ws.Cells.Item(11, 8).Value = -0.0030543261673301

$ws.Cells.Item(12, 1).Value = 1000
$ws.Cells.Item(12, 2).Value = "falling"
$ws.Cells.Item(12, 3).Value = -2.759680032730102
$ws.Cells.Item(12, 4).Value = 9.431608200073242
$ws.Cells.Item(12, 5).Value = 0.08809284307062609
$ws.Cells.Item(12, 6).Value = -0.0395535230636596
$ws.Cells.Item(12, 7).Value = -0.0021380283869802
$ws.Cells.Item(12, 8).Value = 0.0117591563612222

# Remove the old final row (now shifted down to row 32) so the table keeps 30 data rows
$ws.Rows.Item(32).Delete()

# Re-sequence the timestamp column (A) for every data row, 0, 100, 200 ... 2900
For ($r = 2; $r -le 31; $r++) {
    $ws.Cells.Item($r, 1).Value = ($r - 2) * 100
}

Write-Host "done"